# Generate Report for Handoff
# Adds two new rows (for 2ee3cda3-d773-46ab-827f-0535d2cb32a2.md and
# ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md) to the Overview, zh-cn and de-de
# tables/sheets, growing every table range from *1:*3 to *1:*5.

$wb = $excel.ActiveWorkbook

$HYPERLINK_UNDERLINE = 2        # xlUnderlineStyleSingle
$HYPERLINK_COLOR     = 15570276 # BGR for FF6495ED (cornflower blue)
$DATE_FMT            = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "2ee3cda3-d773-46ab-827f-0535d2cb32a2.md"
$wsOverview.Range("B4").Value = "e2e\2ee3cda3-d773-46ab-827f-0535d2cb32a2.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-09-05 02:45:47"
$wsOverview.Range("G4").NumberFormat = $DATE_FMT
$wsOverview.Range("B4").Font.Underline = $HYPERLINK_UNDERLINE
$wsOverview.Range("B4").Font.Color = $HYPERLINK_COLOR

$wsOverview.Range("A5").Value = "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md"
$wsOverview.Range("B5").Value = "e2e\ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-09-05 02:45:47"
$wsOverview.Range("G5").NumberFormat = $DATE_FMT
$wsOverview.Range("B5").Font.Underline = $HYPERLINK_UNDERLINE
$wsOverview.Range("B5").Font.Color = $HYPERLINK_COLOR

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ee3cda3d77346ab827f0535d2cb32a2000000000/e2e/2ee3cda3-d773-46ab-827f-0535d2cb32a2.md", "", "", "e2e\2ee3cda3-d773-46ab-827f-0535d2cb32a2.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea716b3afb984cfeb6d2bb2c07cbc4c5000000000/e2e/ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md", "", "", "e2e\ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = "2ee3cda3-d773-46ab-827f-0535d2cb32a2.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "2ee3cda3-d773-46ab-827f-0535d2cb32a2.d9415d3726a7753d127268b577fe9bc3236ffcef.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-09-05 02:45:42"
$wsZh.Range("H4").NumberFormat = $DATE_FMT
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = $DATE_FMT
$wsZh.Range("M4").Value = "True"
$wsZh.Range("O4").Value = "False"
$wsZh.Range("A4").Font.Underline = $HYPERLINK_UNDERLINE
$wsZh.Range("A4").Font.Color = $HYPERLINK_COLOR

$wsZh.Range("A5").Value = "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.8ede0116b61044dec3730a019f0f767aecc5d1f3.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-09-05 02:45:42"
$wsZh.Range("H5").NumberFormat = $DATE_FMT
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = $DATE_FMT
$wsZh.Range("M5").Value = "True"
$wsZh.Range("O5").Value = "False"
$wsZh.Range("A5").Font.Underline = $HYPERLINK_UNDERLINE
$wsZh.Range("A5").Font.Color = $HYPERLINK_COLOR

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ee3cda3d77346ab827f0535d2cb32a2000000000/e2e/2ee3cda3-d773-46ab-827f-0535d2cb32a2.md", "", "", "2ee3cda3-d773-46ab-827f-0535d2cb32a2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea716b3afb984cfeb6d2bb2c07cbc4c5000000000/e2e/ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md", "", "", "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = "2ee3cda3-d773-46ab-827f-0535d2cb32a2.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "2ee3cda3-d773-46ab-827f-0535d2cb32a2.d9415d3726a7753d127268b577fe9bc3236ffcef.de-de.xlf"
$wsDe.Range("H4").Value = "2016-09-05 02:45:47"
$wsDe.Range("H4").NumberFormat = $DATE_FMT
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = $DATE_FMT
$wsDe.Range("M4").Value = "True"
$wsDe.Range("O4").Value = "False"
$wsDe.Range("A4").Font.Underline = $HYPERLINK_UNDERLINE
$wsDe.Range("A4").Font.Color = $HYPERLINK_COLOR

$wsDe.Range("A5").Value = "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.8ede0116b61044dec3730a019f0f767aecc5d1f3.de-de.xlf"
$wsDe.Range("H5").Value = "2016-09-05 02:45:47"
$wsDe.Range("H5").NumberFormat = $DATE_FMT
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = $DATE_FMT
$wsDe.Range("M5").Value = "True"
$wsDe.Range("O5").Value = "False"
$wsDe.Range("A5").Font.Underline = $HYPERLINK_UNDERLINE
$wsDe.Range("A5").Font.Color = $HYPERLINK_COLOR

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ee3cda3d77346ab827f0535d2cb32a2000000000/e2e/2ee3cda3-d773-46ab-827f-0535d2cb32a2.md", "", "", "2ee3cda3-d773-46ab-827f-0535d2cb32a2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea716b3afb984cfeb6d2bb2c07cbc4c5000000000/e2e/ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md", "", "", "ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md") | Out-Null

Write-Host "Report rows added for 2ee3cda3-d773-46ab-827f-0535d2cb32a2.md and ea716b3a-fb98-4cfe-b6d2-bb2c07cbc4c5.md"
